# Timesheet January 2020 - add the "Mar 2" entries (commit: "Added timesheet of date 02-03-2020")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 236: blank separator row between the "Feb 28" block and the new "Mar 2" block.
# Copy formatting from an existing separator row (A230:C230) so the fill/alignment
# styles (s="5"/"6"/"5") match the rest of the sheet exactly.
$ws.Range("A230:C230").Copy()
$ws.Range("A236:C236").PasteSpecial(-4122)

# Rows 237-244: regular single-line entries. Copy formatting from an existing plain
# data row (A231:C231, styles s="1"/"3"/"1") so no new styles are introduced.
$ws.Range("A231:C231").Copy()
$ws.Range("A237:C244").PasteSpecial(-4122)

# Row 245: final entry of the day - its note wraps onto two lines, so copy formatting
# from an existing wrapped row (A6:C6, styles s="1"/"2"/"1", row height 30).
$ws.Range("A6:C6").Copy()
$ws.Range("A245:C245").PasteSpecial(-4122)

$times = @(
    "Mar 2 10:00 to 11:00",
    "Mar 2 11:00 to 12:00",
    "Mar 2 12:00 to 13:00",
    "Mar 2 13:00 to 14:00",
    "Mar 2 14:00 to 15:00",
    "Mar 2 15:00 to 16:00",
    "Mar 2 16:00 to 17:00",
    "Mar 2 17:00 to 18:00"
)
$tasks = @(
    "Build django backend for predicting automation time, failed units and output units",
    "Working on transforming outfile dates",
    "Discussion of project with frontend developer",
    "Lunch",
    "Issues in transform output unit data, fixing issues",
    "Resolving issue",
    "Resolved issue, working on model building using new modifications",
    "Model building done by modifying some code"
)

$row = 237
for ($i = 0; $i -lt $times.Length; $i++) {
    $ws.Cells.Item($row, 1).Value = $times[$i]
    $ws.Cells.Item($row, 2).Value = $tasks[$i]
    $ws.Cells.Item($row, 3).Value = "Infimetrics"
    $row++
}

# Row 245 - the Task text was entered before the Time text (matches the shared-string
# insertion order captured in the diff), and its height is set explicitly to 30 so the
# wrapped two-line note displays fully.
$ws.Cells.Item(245, 2).Value = "Deployed model successfully, also implemented start time and end time logic`nsuccessfully."
$ws.Cells.Item(245, 1).Value = "Mar 2 18:00 to 19:00"
$ws.Cells.Item(245, 3).Value = "Infimetrics"
$ws.Rows.Item(245).RowHeight = 30

# Restore the view to show the newly-added rows, matching the author's final selection.
$ws.Range("A218").Select()
$ws.Range("B245").Select()
